$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(3, 12).Value = 1.62
$ws.Cells.Item(3, 13).Value = 2.2
$ws.Cells.Item(4, 7).Value = 2.62
$ws.Cells.Item(4, 9).Value = 2.37
$ws.Cells.Item(4, 20).Value = 10.25
$ws.Cells.Item(4, 21).Value = 16
$ws.Cells.Item(4, 22).Value = 10.5
$ws.Cells.Item(4, 23).Value = 32
$ws.Cells.Item(4, 24).Value = 21
$ws.Cells.Item(4, 30).Value = 9.75
$ws.Cells.Item(4, 31).Value = 14
$ws.Cells.Item(4, 33).Value = 26
$ws.Cells.Item(4, 34).Value = 19
$ws.Cells.Item(7, 7).Value = 3
$ws.Cells.Item(7, 9).Value = 2.27
$ws.Cells.Item(7, 16).Value = 1.5
$ws.Cells.Item(7, 17).Value = 2.27
$ws.Cells.Item(7, 20).Value = 7.6
$ws.Cells.Item(7, 22).Value = 11.25
$ws.Cells.Item(7, 27).Value = 6.2
$ws.Cells.Item(7, 33).Value = 22
$ws.Cells.Item(7, 35).Value = 40
$ws.Cells.Item(12, 20).Value = 14
$ws.Cells.Item(12, 21).Value = 21
$ws.Cells.Item(12, 25).Value = 26
$ws.Cells.Item(12, 27).Value = 7.4
$ws.Cells.Item(12, 28).Value = 11.75
$ws.Cells.Item(12, 30).Value = 10
$ws.Cells.Item(12, 31).Value = 11.25
$ws.Cells.Item(12, 33).Value = 18.5
$ws.Cells.Item(12, 35).Value = 20
$ws.Cells.Item(13, 7).Value = 2.22
$ws.Cells.Item(13, 8).Value = 3.35
$ws.Cells.Item(13, 9).Value = 2.9
$ws.Cells.Item(13, 15).Value = 1.88
$ws.Cells.Item(13, 18).Value = 1.6
$ws.Cells.Item(13, 19).Value = 2.07
$ws.Cells.Item(13, 20).Value = 8.5
$ws.Cells.Item(13, 21).Value = 11.5
$ws.Cells.Item(13, 22).Value = 9
$ws.Cells.Item(13, 23).Value = 22
$ws.Cells.Item(13, 24).Value = 17.5
$ws.Cells.Item(13, 25).Value = 25
$ws.Cells.Item(13, 27).Value = 6.6
$ws.Cells.Item(13, 30).Value = 10.75
$ws.Cells.Item(13, 31).Value = 16.5
$ws.Cells.Item(13, 32).Value = 10.25
$ws.Cells.Item(13, 34).Value = 23
$ws.Cells.Item(13, 35).Value = 28
$ws.Cells.Item(13, 36).Value = 350
$ws.Cells.Item(20, 11).Value = 8
$ws.Cells.Item(27, 7).Value = 2.27
$ws.Cells.Item(27, 8).Value = 3
$ws.Cells.Item(27, 9).Value = 3.15
$ws.Cells.Item(27, 12).Value = 1.39
$ws.Cells.Item(27, 13).Value = 2.55
$ws.Cells.Item(27, 14).Value = 2.15
$ws.Cells.Item(27, 15).Value = 1.55
$ws.Cells.Item(27, 16).Value = 1.5
$ws.Cells.Item(27, 17).Value = 2.27
$ws.Cells.Item(27, 20).Value = 6.8
$ws.Cells.Item(27, 21).Value = 10.5
$ws.Cells.Item(27, 23).Value = 23
$ws.Cells.Item(27, 24).Value = 20
$ws.Cells.Item(27, 26).Value = 7.5
$ws.Cells.Item(27, 29).Value = 90
$ws.Cells.Item(27, 30).Value = 8
$ws.Cells.Item(27, 31).Value = 15.5
$ws.Cells.Item(27, 32).Value = 11.5
$ws.Cells.Item(27, 33).Value = 40
$ws.Cells.Item(27, 34).Value = 32
$ws.Cells.Item(28, 9).Value = 3.85
$ws.Cells.Item(28, 14).Value = 1.83
$ws.Cells.Item(28, 15).Value = 1.78
$ws.Cells.Item(28, 18).Value = 1.72
$ws.Cells.Item(28, 19).Value = 1.9
$ws.Cells.Item(28, 20).Value = 7.3
$ws.Cells.Item(28, 25).Value = 26
$ws.Cells.Item(28, 26).Value = 9.75
$ws.Cells.Item(28, 30).Value = 10.75
$ws.Cells.Item(28, 33).Value = 60
$ws.Cells.Item(28, 34).Value = 37
$ws.Cells.Item(29, 7).Value = 2.4
$ws.Cells.Item(29, 8).Value = 3.1
$ws.Cells.Item(29, 9).Value = 2.87
$ws.Cells.Item(29, 12).Value = 1.4
$ws.Cells.Item(29, 13).Value = 2.52
$ws.Cells.Item(29, 14).Value = 2.15
$ws.Cells.Item(29, 15).Value = 1.55
$ws.Cells.Item(29, 16).Value = 1.5
$ws.Cells.Item(29, 17).Value = 2.27
$ws.Cells.Item(29, 18).Value = 1.88
$ws.Cells.Item(29, 19).Value = 1.72
$ws.Cells.Item(29, 20).Value = 6.7
$ws.Cells.Item(29, 21).Value = 10.75
$ws.Cells.Item(29, 22).Value = 9.75
$ws.Cells.Item(29, 23).Value = 25
$ws.Cells.Item(29, 24).Value = 23
$ws.Cells.Item(29, 26).Value = 7.6
$ws.Cells.Item(29, 27).Value = 6
$ws.Cells.Item(29, 28).Value = 16
$ws.Cells.Item(29, 29).Value = 90
$ws.Cells.Item(29, 30).Value = 7.7
$ws.Cells.Item(29, 31).Value = 13.5
$ws.Cells.Item(29, 32).Value = 10.75
$ws.Cells.Item(29, 33).Value = 35
$ws.Cells.Item(29, 34).Value = 27
$ws.Cells.Item(29, 35).Value = 40
$ws.Cells.Item(29, 36).Value = 900
$ws.Cells.Item(31, 10).Value = 1.05
$ws.Cells.Item(31, 12).Value = 1.29
$ws.Cells.Item(32, 19).Value = 2.15
$ws.Cells.Item(32, 21).Value = 15
$ws.Cells.Item(32, 25).Value = 28
$ws.Cells.Item(32, 31).Value = 12.5
$ws.Cells.Item(32, 33).Value = 24
$ws.Cells.Item(32, 34).Value = 17
$ws.Cells.Item(36, 8).Value = 3.25
$ws.Cells.Item(36, 9).Value = 3.2
$ws.Cells.Item(36, 10).Value = 1.06
$ws.Cells.Item(36, 11).Value = 7.1
$ws.Cells.Item(36, 12).Value = 1.32
$ws.Cells.Item(36, 13).Value = 3.15
$ws.Cells.Item(36, 14).Value = 1.93
$ws.Cells.Item(36, 15).Value = 1.78
$ws.Cells.Item(36, 16).Value = 1.44
$ws.Cells.Item(36, 17).Value = 2.62
$ws.Cells.Item(36, 18).Value = 1.75
$ws.Cells.Item(36, 19).Value = 1.95
$ws.Cells.Item(36, 20).Value = 7.5
$ws.Cells.Item(36, 21).Value = 10
$ws.Cells.Item(36, 24).Value = 17
$ws.Cells.Item(36, 25).Value = 28
$ws.Cells.Item(36, 26).Value = 7.1
$ws.Cells.Item(36, 27).Value = 6.4
$ws.Cells.Item(36, 28).Value = 14
$ws.Cells.Item(36, 29).Value = 65
$ws.Cells.Item(36, 30).Value = 9.75
$ws.Cells.Item(36, 31).Value = 17
$ws.Cells.Item(36, 32).Value = 11.25
$ws.Cells.Item(36, 34).Value = 29
$ws.Cells.Item(36, 35).Value = 37
$ws.Cells.Item(36, 36).Value = 500
$ws.Cells.Item(38, 8).Value = 3.35
$ws.Cells.Item(38, 11).Value = 7.3
$ws.Cells.Item(38, 12).Value = 1.31
$ws.Cells.Item(38, 19).Value = 1.98
$ws.Cells.Item(38, 21).Value = 11.25
$ws.Cells.Item(38, 23).Value = 23
$ws.Cells.Item(38, 24).Value = 19
$ws.Cells.Item(38, 26).Value = 7.3
$ws.Cells.Item(38, 27).Value = 6.5
$ws.Cells.Item(38, 28).Value = 14
$ws.Cells.Item(38, 30).Value = 9.25
$ws.Cells.Item(38, 31).Value = 14.5
$ws.Cells.Item(38, 34).Value = 23
$ws.Cells.Item(39, 7).Value = 3.1
$ws.Cells.Item(39, 9).Value = 2.1
$ws.Cells.Item(39, 18).Value = 1.82
$ws.Cells.Item(39, 20).Value = 9.25
$ws.Cells.Item(39, 21).Value = 15.5
$ws.Cells.Item(39, 22).Value = 11.25
$ws.Cells.Item(39, 23).Value = 40
$ws.Cells.Item(39, 24).Value = 28
$ws.Cells.Item(39, 30).Value = 7.3
$ws.Cells.Item(39, 31).Value = 9.75
$ws.Cells.Item(39, 33).Value = 19
$ws.Cells.Item(39, 34).Value = 17.5
$ws.Cells.Item(42, 7).Value = 2.15
$ws.Cells.Item(42, 8).Value = 3.3
$ws.Cells.Item(42, 9).Value = 3.1
$ws.Cells.Item(42, 13).Value = 2.65
$ws.Cells.Item(42, 21).Value = 9.5
$ws.Cells.Item(42, 23).Value = 20
$ws.Cells.Item(42, 24).Value = 19
$ws.Cells.Item(42, 26).Value = 8.25
$ws.Cells.Item(42, 27).Value = 6.4
$ws.Cells.Item(42, 28).Value = 16.5
$ws.Cells.Item(42, 29).Value = 90
$ws.Cells.Item(42, 32).Value = 11.5
$ws.Cells.Item(42, 34).Value = 30
$ws.Cells.Item(42, 35).Value = 45
